$d = $word.ActiveDocument

$d.Content.Find.Execute("49×67=3283", $true, $false, $false, $false, $false, $true, 1, $false, "27×12=324", 2) | Out-Null
$d.Content.Find.Execute("89×86=7654", $true, $false, $false, $false, $false, $true, 1, $false, "15×38=570", 2) | Out-Null
$d.Content.Find.Execute("35×86=3010", $true, $false, $false, $false, $false, $true, 1, $false, "21×20=420", 2) | Out-Null
$d.Content.Find.Execute("36×45=1620", $true, $false, $false, $false, $false, $true, 1, $false, "37×27=999", 2) | Out-Null
$d.Content.Find.Execute("34×25=850", $true, $false, $false, $false, $false, $true, 1, $false, "38×59=2242", 2) | Out-Null
$d.Content.Find.Execute("89×79=7031", $true, $false, $false, $false, $false, $true, 1, $false, "62×64=3968", 2) | Out-Null
$d.Content.Find.Execute("53×30=1590", $true, $false, $false, $false, $false, $true, 1, $false, "11×30=330", 2) | Out-Null
$d.Content.Find.Execute("83×73=6059", $true, $false, $false, $false, $false, $true, 1, $false, "31×75=2325", 2) | Out-Null
$d.Content.Find.Execute("17×21=357", $true, $false, $false, $false, $false, $true, 1, $false, "65×89=5785", 2) | Out-Null
$d.Content.Find.Execute("21×80=1680", $true, $false, $false, $false, $false, $true, 1, $false, "78×33=2574", 2) | Out-Null
$d.Content.Find.Execute("76×39=2964", $true, $false, $false, $false, $false, $true, 1, $false, "32×32=1024", 2) | Out-Null
$d.Content.Find.Execute("65×25=1625", $true, $false, $false, $false, $false, $true, 1, $false, "48×53=2544", 2) | Out-Null
$d.Content.Find.Execute("20×56=1120", $true, $false, $false, $false, $false, $true, 1, $false, "61×39=2379", 2) | Out-Null
$d.Content.Find.Execute("96×38=3648", $true, $false, $false, $false, $false, $true, 1, $false, "87×37=3219", 2) | Out-Null
$d.Content.Find.Execute("88×12=1056", $true, $false, $false, $false, $false, $true, 1, $false, "20×25=500", 2) | Out-Null
$d.Content.Find.Execute("97×35=3395", $true, $false, $false, $false, $false, $true, 1, $false, "36×20=720", 2) | Out-Null
$d.Content.Find.Execute("31×40=1240", $true, $false, $false, $false, $false, $true, 1, $false, "41×38=1558", 2) | Out-Null
$d.Content.Find.Execute("64×99=6336", $true, $false, $false, $false, $false, $true, 1, $false, "38×21=798", 2) | Out-Null
$d.Content.Find.Execute("60×79=4740", $true, $false, $false, $false, $false, $true, 1, $false, "45×40=1800", 2) | Out-Null
$d.Content.Find.Execute("95×86=8170", $true, $false, $false, $false, $false, $true, 1, $false, "66×40=2640", 2) | Out-Null
$d.Content.Find.Execute("91×16=1456", $true, $false, $false, $false, $false, $true, 1, $false, "89×76=6764", 2) | Out-Null
$d.Content.Find.Execute("54×66=3564", $true, $false, $false, $false, $false, $true, 1, $false, "14×57=798", 2) | Out-Null
$d.Content.Find.Execute("83×70=5810", $true, $false, $false, $false, $false, $true, 1, $false, "87×21=1827", 2) | Out-Null
$d.Content.Find.Execute("47×74=3478", $true, $false, $false, $false, $false, $true, 1, $false, "39×95=3705", 2) | Out-Null
$d.Content.Find.Execute("91×38=3458", $true, $false, $false, $false, $false, $true, 1, $false, "63×49=3087", 2) | Out-Null
